$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Technique labels (column A) ---
# "Marching Cubes" (rows 2-13)   -> "Edge Detection" (data unchanged)
# "Watershed"      (rows 14-25)  -> "Otsu" (new data for rows 17-19, 23-25)
# "Edge Detection" (rows 26-37)  -> "Watershed" (gets old "Marching Cubes" data)
$ws.Range("A2:A13").Value2  = "Edge Detection"
$ws.Range("A14:A25").Value2 = "Otsu"
$ws.Range("A26:A37").Value2 = "Watershed"

# --- Numeric data (columns D, E) ---

# Rows 2-13 now hold what used to be the rows 26-37 values
$ws.Range("D2").Value2  = 14.533595111867021
$ws.Range("E2").Value2  = 16.148542165749522
$ws.Range("D3").Value2  = 102.7973195300453
$ws.Range("E3").Value2  = 99.284519657733725
$ws.Range("D4").Value2  = 101.8230655722588
$ws.Range("E4").Value2  = 101.4611836238015
$ws.Range("D5").Value2  = 16.813198102189322
$ws.Range("E5").Value2  = 20.0376588623262
$ws.Range("D6").Value2  = 100.2943306115923
$ws.Range("E6").Value2  = 108.1059908672057
$ws.Range("D7").Value2  = 50.176953938103843
$ws.Range("E7").Value2  = 65.902816281114156
$ws.Range("D8").Value2  = 13.84499381841573
$ws.Range("E8").Value2  = 16.212580140341039
$ws.Range("D9").Value2  = 150.77097147398121
$ws.Range("E9").Value2  = 157.02951863135701
$ws.Range("D10").Value2 = 62.187452474357329
$ws.Range("E10").Value2 = 82.478070129762301
$ws.Range("D11").Value2 = 13.012877906533349
$ws.Range("E11").Value2 = 16.23338314328134
$ws.Range("D12").Value2 = 101.59066637534571
$ws.Range("E12").Value2 = 72.585021128572961
$ws.Range("D13").Value2 = 57.21813766268955
$ws.Range("E13").Value2 = 48.36376403274064

# Rows 17-19 ("Otsu" / Heart failure without infarct) get new values
$ws.Range("D17").Value2 = 32.867396690934299
$ws.Range("E17").Value2 = 18.87231263031746
$ws.Range("D18").Value2 = 428.17025270436721
$ws.Range("E18").Value2 = 482.57526238321861
$ws.Range("D19").Value2 = 102.9561926671045
$ws.Range("E19").Value2 = 132.2540520011988

# Rows 23-25 ("Otsu" / Normal) get new values
$ws.Range("D23").Value2 = 17.861086752801722
$ws.Range("E23").Value2 = 18.126721849355171
$ws.Range("D24").Value2 = 293.86695966782639
$ws.Range("E24").Value2 = 298.87144095546188
$ws.Range("D25").Value2 = 56.02734316913179
$ws.Range("E25").Value2 = 44.214879453982093

# Rows 26-37 now hold what used to be the rows 2-13 values
$ws.Range("D26").Value2 = 12.03768911223059
$ws.Range("E26").Value2 = 13.53116963868135
$ws.Range("D27").Value2 = 195.98449100186849
$ws.Range("E27").Value2 = 16.23740465086442
$ws.Range("D28").Value2 = 139.12969927417731
$ws.Range("E28").Value2 = 14.880036838303241
$ws.Range("D29").Value2 = 15.110845585629679
$ws.Range("E29").Value2 = 16.791451045920741
$ws.Range("D30").Value2 = 165.8681147955904
$ws.Range("E30").Value2 = 27.278996653032319
$ws.Range("D31").Value2 = 117.05796284591619
$ws.Range("E31").Value2 = 18.51043123435031
$ws.Range("D32").Value2 = 19.655052036391329
$ws.Range("E32").Value2 = 20.12392812487505
$ws.Range("D33").Value2 = 54.748522663978591
$ws.Range("E33").Value2 = 20.573987391585149
$ws.Range("D34").Value2 = 13.90189030222141
$ws.Range("E34").Value2 = 12.89567834873691
$ws.Range("D35").Value2 = 22.0318377319491
$ws.Range("E35").Value2 = 15.85739820227772
$ws.Range("D36").Value2 = 44.987293088577673
$ws.Range("E36").Value2 = 23.73204003206974
$ws.Range("D37").Value2 = 14.033924547349899
$ws.Range("E37").Value2 = 20.009258906245201

# Clear the lingering selection so the sheet view defaults back to A1,
# and let column A's width recompute (best fit) for the new, shorter
# technique name ("Edge Detection" replacing "Marching Cubes").
$ws.Range("A1").Select()
$ws.Columns("A").AutoFit()
